$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reclassify the child-table header row: positional placeholders become
# named-column placeholders (table column classification for Excel output).
$ws.Range("A8").Value = "{child[ID]:linked}"
$ws.Range("B8").Value = "{child[DESCRIPTION]:linked}"
$ws.Range("C8").Value = "{child[DATE]}"

# Move the active selection from P2:P3 to D8 (also clears the stale
# horizontal scroll position at D1).
$ws.Range("D8").Select()
